$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Main data edit: update the "days" parameter used by the WORKDAY() formulas
# in columns AK/AO/AP/AS (rows 8-28). Changing this single input cell drives
# the recalculation of all dependent date/formula cells.
$ws.Range("AM2").Value = 20

# Reflect the cell the user ended up with selected after making the change.
$ws.Range("AK9").Select()
